$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.243.29"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "2.014.29"
$ws.Range("E3").Value = "  -1.26%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.28"
$ws.Range("E5").Value = "  +2.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.641"
$ws.Range("E6").Value = "  -3.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.84"
$ws.Range("E7").Value = "  +11.99%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "59.17"
$ws.Range("E9").Value = "  -6.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.369"
$ws.Range("E10").Value = "  +0.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0745"
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.104"
$ws.Range("E12").Value = "  -1.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.906"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.76"
$ws.Range("E14").Value = "  +3.89%  "
$ws.Range("D15").Value = "2.306.39"
$ws.Range("E15").Value = "  -1.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.43"
$ws.Range("E16").Value = "  +0.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.84"
$ws.Range("E17").Value = "  +12.99%  "
$ws.Range("D18").Value = "2.050.65"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").Value = "36.167.83"
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.06"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("D21").Value = "0.0₃0858"
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("E22").Value = "  +1.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.88"
$ws.Range("E23").Value = "  -1.14%  "
$ws.Range("E24").Value = "  +19.44%  "
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("E26").Value = "  -2.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.55"
$ws.Range("E27").Value = "  +3.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.99"
$ws.Range("E28").Value = "  -0.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.62"
$ws.Range("E29").Value = "  -1.54%  "
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.112"
$ws.Range("E31").Value = "  +27.23%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.13"
$ws.Range("E32").Value = "  +3.31%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.18"
$ws.Range("E33").Value = "  -2.01%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.52"
$ws.Range("E34").Value = "  +15.09%  "
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.53"
$ws.Range("E35").Value = "  +2.99%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0604"
$ws.Range("E36").Value = "  +1.26%  "
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.80"
$ws.Range("E38").Value = "  -0.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.86"
$ws.Range("E39").Value = "  +15.59%  "
$ws.Range("E40").Value = "  +13.42%  "
$ws.Range("E41").Value = "  +0.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.91"
$ws.Range("E42").Value = "  +0.87%  "
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.13"
$ws.Range("E44").Value = "  +2.54%  "
$ws.Range("E45").Value = "  +4.84%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.86"
$ws.Range("E46").Value = "  +6.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "94.05"
$ws.Range("E47").Value = "  +0.73%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "1.423.93"
$ws.Range("E48").Value = "  +3.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.58"
$ws.Range("E49").Value = "  +14.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.92"
$ws.Range("E50").Value = "  -0.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.41"
$ws.Range("E51").Value = "  +3.37%  "
